$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("India Super League")

# ---------------------------------------------------------------------
# Row 116 (existing row, id=114) gets its data replaced with a new match
# (Mohun Bagan SG vs Chennaiyin FC) plus extra columns H/I/J now filled.
# ---------------------------------------------------------------------
$ws.Cells.Item(116, 2).Value  = 7749761                     # B id
$ws.Cells.Item(116, 3).Value  = "India Super League"         # C Div
$ws.Cells.Item(116, 4).Value  = "India Super League"         # D Div Original Name
$ws.Cells.Item(116, 5).Value  = 45382.45833333334            # E Date
$ws.Cells.Item(116, 6).Value  = "Mohun Bagan SG"              # F HomeTeam
$ws.Cells.Item(116, 7).Value  = "Chennaiyin FC"               # G AwayTeam
$ws.Cells.Item(116, 8).Value  = 2                             # H FTHG
$ws.Cells.Item(116, 9).Value  = 3                             # I FTAG
$ws.Cells.Item(116, 10).Value = "A"                            # J FTR
$ws.Cells.Item(116, 11).Value = 1.533                          # K oddH_op
$ws.Cells.Item(116, 12).Value = 4                              # L oddD_op
$ws.Cells.Item(116, 13).Value = 5.5                            # M oddA_op
$ws.Cells.Item(116, 14).Value = 1.45                           # N oddH
$ws.Cells.Item(116, 15).Value = 4.333                          # O oddD
$ws.Cells.Item(116, 16).Value = 5.75                           # P oddA
$ws.Cells.Item(116, 17).Value = -1.25                          # Q Ah
$ws.Cells.Item(116, 18).Value = 2.025                          # R oddAHH
$ws.Cells.Item(116, 19).Value = 1.825                          # S oddAHA
$ws.Cells.Item(116, 20).Value = 3                              # T AhOU
$ws.Cells.Item(116, 21).Value = 1.925                          # U oddAHOver
$ws.Cells.Item(116, 22).Value = 1.925                          # V oddAHUnder
$ws.Cells.Item(116, 23).Value = -1                             # W PLH
$ws.Cells.Item(116, 24).Value = -1                             # X PLD
$ws.Cells.Item(116, 25).Value = 4.75                           # Y PLA
$ws.Cells.Item(116, 26).Value = -1                             # Z PL_Ahh
$ws.Cells.Item(116, 27).Value = 0.825                          # AA PL_Aha
$ws.Cells.Item(116, 28).Value = 0.925                          # AB PL_AhOver
$ws.Cells.Item(116, 29).Value = -1                             # AC PL_AhUnder

# ---------------------------------------------------------------------
# Row 117 (new, id=115) Hyderabad FC vs Mumbai City FC
# ---------------------------------------------------------------------
$ws.Cells.Item(117, 1).Value  = 115                            # A id (row index)
$ws.Cells.Item(117, 2).Value  = 7749875                        # B id
$ws.Cells.Item(117, 3).Value  = "India Super League"           # C Div
$ws.Cells.Item(117, 4).Value  = "India Super League"           # D Div Original Name
$ws.Cells.Item(117, 5).Value  = 45383.45833333334              # E Date
$ws.Cells.Item(117, 6).Value  = "Hyderabad FC"                 # F HomeTeam
$ws.Cells.Item(117, 7).Value  = "Mumbai City FC"                # G AwayTeam
$ws.Cells.Item(117, 8).Value  = 0                               # H FTHG
$ws.Cells.Item(117, 9).Value  = 3                               # I FTAG
$ws.Cells.Item(117, 10).Value = "A"                              # J FTR
$ws.Cells.Item(117, 11).Value = 9.5                              # K oddH_op
$ws.Cells.Item(117, 12).Value = 5.5                              # L oddD_op
$ws.Cells.Item(117, 13).Value = 1.25                             # M oddA_op
$ws.Cells.Item(117, 14).Value = 7                                # N oddH
$ws.Cells.Item(117, 15).Value = 5.5                              # O oddD
$ws.Cells.Item(117, 16).Value = 1.333                            # P oddA
$ws.Cells.Item(117, 17).Value = 1.5                              # Q Ah
$ws.Cells.Item(117, 18).Value = 1.85                             # R oddAHH
$ws.Cells.Item(117, 19).Value = 1.95                             # S oddAHA
$ws.Cells.Item(117, 20).Value = 3.25                             # T AhOU
$ws.Cells.Item(117, 21).Value = 2                                # U oddAHOver
$ws.Cells.Item(117, 22).Value = 1.8                              # V oddAHUnder
$ws.Cells.Item(117, 23).Value = -1                               # W PLH
$ws.Cells.Item(117, 24).Value = -1                               # X PLD
$ws.Cells.Item(117, 25).Value = 0.333                            # Y PLA
$ws.Cells.Item(117, 26).Value = -1                               # Z PL_Ahh
$ws.Cells.Item(117, 27).Value = 0.95                             # AA PL_Aha
$ws.Cells.Item(117, 28).Value = -0.5                             # AB PL_AhOver
$ws.Cells.Item(117, 29).Value = 0.4                              # AC PL_AhUnder

# ---------------------------------------------------------------------
# Row 118 (new, id=116) Odisha FC vs Punjab FC  (not played yet - no
# FTHG/FTAG/FTR)
# ---------------------------------------------------------------------
$ws.Cells.Item(118, 1).Value  = 116                             # A id (row index)
$ws.Cells.Item(118, 2).Value  = 7749773                         # B id
$ws.Cells.Item(118, 3).Value  = "India Super League"            # C Div
$ws.Cells.Item(118, 4).Value  = "India Super League"            # D Div Original Name
$ws.Cells.Item(118, 5).Value  = 45384.45833333334               # E Date
$ws.Cells.Item(118, 6).Value  = "Odisha FC"                     # F HomeTeam
$ws.Cells.Item(118, 7).Value  = "Punjab FC"                     # G AwayTeam
$ws.Cells.Item(118, 11).Value = 1.65                            # K oddH_op
$ws.Cells.Item(118, 12).Value = 3.75                            # L oddD_op
$ws.Cells.Item(118, 13).Value = 4.2                             # M oddA_op
$ws.Cells.Item(118, 14).Value = 1.666                           # N oddH
$ws.Cells.Item(118, 15).Value = 3.75                            # O oddD
$ws.Cells.Item(118, 16).Value = 5                               # P oddA
$ws.Cells.Item(118, 17).Value = -0.75                           # Q Ah
$ws.Cells.Item(118, 18).Value = 1.85                            # R oddAHH
$ws.Cells.Item(118, 19).Value = 1.95                            # S oddAHA
$ws.Cells.Item(118, 20).Value = 2.75                            # T AhOU
$ws.Cells.Item(118, 21).Value = 1.875                           # U oddAHOver
$ws.Cells.Item(118, 22).Value = 1.925                           # V oddAHUnder
$ws.Cells.Item(118, 23).Value = 0                               # W PLH
$ws.Cells.Item(118, 24).Value = 0                               # X PLD
$ws.Cells.Item(118, 25).Value = 0                               # Y PLA
$ws.Cells.Item(118, 26).Value = 0                               # Z PL_Ahh
$ws.Cells.Item(118, 27).Value = 0                               # AA PL_Aha

# Apply the same formatting as the other "id" column cells (bold/border)
# and date format used throughout the sheet, matching existing rows.
$ws.Cells.Item(115, 1).Copy()
$ws.Cells.Item(117, 1).PasteSpecial(-4122)
$ws.Cells.Item(118, 1).PasteSpecial(-4122)

$ws.Cells.Item(115, 5).Copy()
$ws.Cells.Item(117, 5).PasteSpecial(-4122)
$ws.Cells.Item(118, 5).PasteSpecial(-4122)

$excel.CutCopyMode = 0
